# Auto-generated edit script: updates cryptocurrency price/volume table
# to reflect the latest GitHub Actions scrape (Sun Jul 28 17:07:29 UTC 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.800.51"
$ws.Range("E2").Value = "  -1.34%  "

$ws.Range("D3").Value = "3.270.89"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.10"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.14"
$ws.Range("E6").Value = "  +1.06%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.97%  "

$ws.Range("E9").Value = "  -2.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.58"
$ws.Range("E10").Value = "  -1.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.409"
$ws.Range("E11").Value = "  -4.01%  "

$ws.Range("D12").Value = "3.834.63"
$ws.Range("E12").Value = "  -0.28%  "

$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.41"
$ws.Range("E14").Value = "  -4.41%  "

$ws.Range("D15").Value = "67.777.70"
$ws.Range("E15").Value = "  -1.28%  "

$ws.Range("E16").Value = "  -1.85%  "

$ws.Range("D17").Value = "3.254.74"
$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.71"
$ws.Range("E18").Value = "  -2.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.41"
$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "402.47"
$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.55"
$ws.Range("E21").Value = "  -1.83%  "

$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.15"
$ws.Range("E23").Value = "  -1.18%  "

$ws.Range("E24").Value = "  -1.63%  "

$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.188"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.50"
$ws.Range("E27").Value = "  -1.60%  "

$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.94"
$ws.Range("E29").Value = "  -1.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.71"
$ws.Range("E30").Value = "  -0.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.48"
$ws.Range("E31").Value = "  -3.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.91"
$ws.Range("E32").Value = "  -3.00%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("E34").Value = "  -3.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "164.24"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("E36").Value = "  -3.65%  "

$ws.Range("E37").Value = "  -3.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.03"
$ws.Range("E38").Value = "  +2.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.803"
$ws.Range("E39").Value = "  -3.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.50"
$ws.Range("E40").Value = "  -2.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.33"
$ws.Range("E41").Value = "  -3.11%  "

$ws.Range("D42").Value = "2.680.21"
$ws.Range("E42").Value = "  +2.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.82"
$ws.Range("E43").Value = "  -1.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("E44").Value = "  -3.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0678"
$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "335.42"
$ws.Range("E46").Value = "  -2.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.62"
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0274"
$ws.Range("E48").Value = "  -2.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.30"
$ws.Range("E49").Value = "  -0.38%  "

$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.967"
$ws.Range("E51").Value = "  -1.42%  "
